$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 417
$ws1.Range("F5").Value = 1244
$ws1.Range("F7").Value = 7504
$ws1.Range("F11").Value = 8171
$ws1.Range("F14").Value = 5578
$ws1.Range("F16").Value = 2541
$ws1.Range("F17").Value = 1091
$ws1.Range("F19").Value = 318
$ws1.Range("F23").Value = 456
$ws1.Range("F24").Value = 1826
$ws1.Range("F25").Value = 31
$ws1.Range("F26").Value = 2705
$ws1.Range("F28").Value = 309
$ws1.Range("F30").Value = 250
$ws1.Range("F31").Value = 622
$ws1.Range("F33").Value = 529
$ws1.Range("F34").Value = 1605
$ws1.Range("F37").Value = 2552

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 100
$ws2.Range("F4").Value = 31

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 417
$ws4.Range("F7").Value = 1244
$ws4.Range("F9").Value = 7504
$ws4.Range("F13").Value = 8171
$ws4.Range("F16").Value = 5578
$ws4.Range("F18").Value = 2541
$ws4.Range("F19").Value = 1091
$ws4.Range("F21").Value = 318
$ws4.Range("F26").Value = 100
$ws4.Range("F27").Value = 456
$ws4.Range("F28").Value = 1826
$ws4.Range("F29").Value = 31
$ws4.Range("F30").Value = 2705
$ws4.Range("F32").Value = 309
$ws4.Range("F34").Value = 250
$ws4.Range("F35").Value = 31
$ws4.Range("F36").Value = 622
$ws4.Range("F38").Value = 529
$ws4.Range("F40").Value = 1605
$ws4.Range("F43").Value = 2552
